$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the apostrophe typo in "Pend d'Oreille River" (row 53, column A)
$ws.Range("A53").Value = "Pend dOreille River"

# Remove the two duplicate "roving station" rows that were accidentally
# entered twice (once with an apostrophe in the name, once without), and
# give every subsequent roving station its own correct Region / Agency
# label by removing the stray duplicate rows so everything shifts up into
# its proper place.
# Row 72 = "Surveyor's Lake" (duplicate of "Surveyors Lake")
$ws.Rows.Item(72).Delete()
# Row 68 = "St. Mary's Lake" (duplicate of "St. Marys Lake")
$ws.Rows.Item(68).Delete()
